# Update "想去人数" (want-to-go count) figures in both the "展览" sheet
# and the merged "全部类型" sheet to reflect the latest scrape.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 16807
$wsExpo.Range("F3").Value = 357
$wsExpo.Range("F4").Value = 748
$wsExpo.Range("F5").Value = 257
$wsExpo.Range("F6").Value = 735
$wsExpo.Range("F7").Value = 1811

# --- Sheet "全部类型" (All Types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 16807
$wsAll.Range("F3").Value = 357
$wsAll.Range("F4").Value = 748
$wsAll.Range("F5").Value = 257
$wsAll.Range("F8").Value = 735
$wsAll.Range("F9").Value = 1811
